# Updated run for publication: refresh the frequency-table values on Sheet1
# (rows 2-5, columns B-X) to the newly recomputed base-frequency figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.0161662817551963
$ws.Range("C2").Value = 0.0023094688221709
$ws.Range("E2").Value = 0.0161662817551963
$ws.Range("F2").Value = 0.0023094688221709
$ws.Range("K2").Value = 0.0138568129330254
$ws.Range("L2").Value = 0.0161662817551963
$ws.Range("M2").Value = 0.0046189376443418
$ws.Range("N2").Value = 0.0046189376443418
$ws.Range("O2").Value = 0.344110854503464
$ws.Range("S2").Value = 0.979214780600462
$ws.Range("T2").Value = 0.0161662817551963
$ws.Range("U2").Value = 0.187066974595843
$ws.Range("V2").Value = 0.979214780600462
$ws.Range("W2").Value = 0.0023094688221709
$ws.Range("X2").Value = 0.97459584295612
$ws.Range("B3").Value = 0.889145496535797
$ws.Range("C3").Value = 0.993071593533487
$ws.Range("D3").Value = 0.0207852193995381
$ws.Range("E3").Value = 0.0115473441108545
$ws.Range("F3").Value = 0.0023094688221709
$ws.Range("H3").Value = 0.092378752886836
$ws.Range("I3").Value = 0.826789838337182
$ws.Range("L3").Value = 0.0046189376443418
$ws.Range("N3").Value = 0.0161662817551963
$ws.Range("O3").Value = 0.0161662817551963
$ws.Range("P3").Value = 0.0092378752886836
$ws.Range("Q3").Value = 0.89838337182448
$ws.Range("R3").Value = 0.0046189376443418
$ws.Range("T3").Value = 0.905311778290993
$ws.Range("U3").Value = 0.79445727482679
$ws.Range("V3").Value = 0.0184757505773672
$ws.Range("W3").Value = 0.909930715935335
$ws.Range("X3").Value = 0.0161662817551963
$ws.Range("C4").Value = 0.0023094688221709
$ws.Range("D4").Value = 0.0023094688221709
$ws.Range("F4").Value = 0.995381062355658
$ws.Range("K4").Value = 0.981524249422633
$ws.Range("M4").Value = 0.993071593533487
$ws.Range("N4").Value = 0.979214780600462
$ws.Range("O4").Value = 0.635103926096998
$ws.Range("P4").Value = 0.0023094688221709
$ws.Range("Q4").Value = 0.0184757505773672
$ws.Range("R4").Value = 0.0207852193995381
$ws.Range("S4").Value = 0.0184757505773672
$ws.Range("T4").Value = 0.0023094688221709
$ws.Range("U4").Value = 0.0161662817551963
$ws.Range("X4").Value = 0.0046189376443418
$ws.Range("B5").Value = 0.0946882217090069
$ws.Range("C5").Value = 0.0023094688221709
$ws.Range("D5").Value = 0.976905311778291
$ws.Range("E5").Value = 0.972286374133949
$ws.Range("H5").Value = 0.907621247113164
$ws.Range("I5").Value = 0.173210161662818
$ws.Range("K5").Value = 0.0046189376443418
$ws.Range("L5").Value = 0.979214780600462
$ws.Range("M5").Value = 0.0023094688221709
$ws.Range("O5").Value = 0.0046189376443418
$ws.Range("P5").Value = 0.988452655889146
$ws.Range("Q5").Value = 0.0831408775981524
$ws.Range("R5").Value = 0.97459584295612
$ws.Range("S5").Value = 0.0023094688221709
$ws.Range("T5").Value = 0.0739030023094688
$ws.Range("U5").Value = 0.0023094688221709
$ws.Range("V5").Value = 0.0023094688221709
$ws.Range("W5").Value = 0.0877598152424942
$ws.Range("X5").Value = 0.0046189376443418
